$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old header row (A1:I1) contents - it's no longer part of the data.
# Rows 2 and 3 (E2:H3) stay exactly where they are.
$ws.Range("A1:I1").ClearContents()

# Append a new data block in column C (rows 27-39), simulating pandas
# reading a dataframe that was written starting at a certain row/col.
$ws.Range("C27").Value = "uF/cm2"
$ws.Range("C28").Value = 0.60454094899999999
$ws.Range("C29").Value = 0.62367113500000004
$ws.Range("C30").Value = 0.61469258000000004
$ws.Range("C31").Value = 0.605172606
$ws.Range("C32").Value = 0.60945885
$ws.Range("C33").Value = 0.60697734000000003
$ws.Range("C34").Value = 0.62096403300000003
$ws.Range("C35").Value = 0.616542433
$ws.Range("C36").Value = 0.61284272799999995
$ws.Range("C37").Value = 0.60666151199999996
$ws.Range("C38").Value = "jarrett"
$ws.Range("C39").Value = "goh"

# Update the view to match the scrolled/selected position after the edit
$ws.Range("E37").Select()
$excel.ActiveWindow.ScrollRow = 10
